$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.428.99'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.382.30'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.88%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.99'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.79%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.49'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.91%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.382.75'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.80%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.66%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.49'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.90%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.125'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.389'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.29%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.960.16'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.97%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.46%  '

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.63%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.384.39'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.21%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.80'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.63%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.548.12'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.08'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.36%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.86'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.36'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '375.62'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.556'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.45%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.524.77'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.08%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +8.18%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '71.21'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.17%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.16%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.49'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.12%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.02%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +5.38%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.79%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.93%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.45'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.32%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.26'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.37%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.54'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.82'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '165.54'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.74%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0776'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.97%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.08%  '

$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.776'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.74%  '

$ws.Range("B43").Value = 'ONDO'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.22'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.29%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +8.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.41'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.71%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.38'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.34%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.75'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +8.99%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.82'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.74%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.74'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.342.36'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +5.14%  '

$ws.Range("B51").Value = 'LidoDAOToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.37'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.37%  '
